# IT Risk Analysis Sheet edits
# 1. Update the risk-explanation text for the "Categories" risk item (I7)
# 2. Move the "Timeline" risk item from the Implementation Risk group (old row 20)
#    up into the Organizational Complexity group (row 18), and remove it from
#    its old spot -- i.e. "Changed time [Timeline] to organizational risk".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Shared-string text edit -----------------------------------------
$ws.Range("I7").Value2 = "The team is unable to provide the correct tools to upload profiles and CTF challenges"
$ws.Rows.Item(7).RowHeight = 91

# --- 2. Move the Timeline row (old row 20) into row 18 -------------------
# Row 18 is currently an empty spacer row; row 20 holds the Timeline data.
# Clear row 18 first, then cut row 20's cells (values + formatting) into it.
$ws.Range("A18:I18").Clear()
$ws.Range("A20:I20").Cut($ws.Range("A18:I18"))

# Cut() does not keep the relative formula for D18/D20, restore it explicitly
$ws.Range("D18").Formula = "=C18*`$B`$3*B18"

# Cut() also does not move the row-level height, only the cell formatting --
# carry that over by hand so row 18 matches how the Timeline row used to look.
$ws.Rows.Item(18).RowHeight = 78

# Row 20 is now vacated -- clear any leftover formatting and collapse its
# height back down so it reads as an ordinary blank row.
$ws.Range("A20:I20").Clear()
$ws.Rows.Item(20).RowHeight = 12.5

# --- 3. Fix up the SUM() formulas that referenced the old layout ---------
# Organizational Complexity total (row 14) now also includes row 18
$ws.Range("D14").Formula = "=SUM(D15:D18)"
# Implementation Risk total (row 19) no longer includes (now-empty) row 20
$ws.Range("D19").Formula = "=SUM(D21)"

# --- 4. Restore the view/selection state left in the saved file ----------
$excel.ActiveWindow.ScrollRow = 18
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("D20").Select()
